$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 471.33334
$ws.Range("I4").Value = 36
$ws.Range("K4").Value = 36
$ws.Range("M4").Value = 78
$ws.Range("H5").Value = 300
$ws.Range("J5").Value = 233.33333
$ws.Range("L5").Value = 233.33333
$ws.Range("N5").Value = -463.33333
$ws.Range("H19").Value = 1415.3
$ws.Range("I19").Value = 200
$ws.Range("K19").Value = 200
$ws.Range("M19").Value = -25
$ws.Range("H32").Value = 3860.4783
$ws.Range("J32").Value = 3899.6
$ws.Range("L32").Value = 3899.6
$ws.Range("N32").Value = -4551.6
$ws.Range("H64").Value = 8308.125
$ws.Range("I64").Value = 6497.5
$ws.Range("J64").Value = 8911.666999999999
$ws.Range("K64").Value = 6497.5
$ws.Range("L64").Value = 8911.666999999999
$ws.Range("M64").Value = -6249.5
$ws.Range("N64").Value = -9407.666999999999
$ws.Range("H67").Value = 8308.125
$ws.Range("I67").Value = 6497.5
$ws.Range("J67").Value = 8911.666999999999
$ws.Range("K67").Value = 6497.5
$ws.Range("L67").Value = 8911.666999999999
$ws.Range("M67").Value = -5639.5
$ws.Range("N67").Value = -10627.667
$ws.Range("H70").Value = 6898.0835
$ws.Range("I70").Value = 1400
$ws.Range("J70").Value = 7997.7
$ws.Range("K70").Value = 4200
$ws.Range("L70").Value = 23993.1
$ws.Range("M70").Value = -3930
$ws.Range("N70").Value = -24533.1
$ws.Range("H73").Value = 6898.0835
$ws.Range("I73").Value = 1400
$ws.Range("J73").Value = 7997.7
$ws.Range("K73").Value = 4200
$ws.Range("L73").Value = 23993.1
$ws.Range("M73").Value = -3264
$ws.Range("N73").Value = -25865.1
$ws.Range("H86").Value = 7196.0557
$ws.Range("I86").Value = 6690.857
$ws.Range("K86").Value = 6690.857
$ws.Range("M86").Value = -5567.857
$ws.Range("H89").Value = 7196.0557
$ws.Range("I89").Value = 6690.857
$ws.Range("K89").Value = 33454.285
$ws.Range("M89").Value = -27838.285
$ws.Range("H138").Value = 3352
$ws.Range("J138").Value = 3558.1785
$ws.Range("L138").Value = 10674.5355
$ws.Range("N138").Value = -20954.5355

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9046.200000000001
$ws.Range("I32").Value = 5166.4253
$ws.Range("K32").Value = 5166.4253
$ws.Range("M32").Value = -4879.4253
$ws.Range("H45").Value = 6853929
$ws.Range("I45").Value = 11989753
$ws.Range("J45").Value = 6164
$ws.Range("K45").Value = 11989753
$ws.Range("L45").Value = 6164
$ws.Range("M45").Value = -11989376
$ws.Range("N45").Value = -6918

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3335102
$ws.Range("I86").Value = 4168448.2
$ws.Range("K86").Value = 4168448.2
$ws.Range("M86").Value = -4167325.2
$ws.Range("H89").Value = 3335102
$ws.Range("I89").Value = 4168448.2
$ws.Range("K89").Value = 20842241
$ws.Range("M89").Value = -20836625
$ws.Range("H134").Value = 10015.889
$ws.Range("I134").Value = 5315
$ws.Range("J134").Value = 12366.333
$ws.Range("K134").Value = 15945
$ws.Range("L134").Value = 37098.999
$ws.Range("M134").Value = -13410
$ws.Range("N134").Value = -42168.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 572.55554
$ws.Range("I7").Value = 389.25
$ws.Range("K7").Value = 389.25
$ws.Range("M7").Value = -276.25
$ws.Range("H16").Value = 2190.8572
$ws.Range("I16").Value = 1172.5
$ws.Range("J16").Value = 2598.2
$ws.Range("K16").Value = 1172.5
$ws.Range("L16").Value = 2598.2
$ws.Range("M16").Value = -885.5
$ws.Range("N16").Value = -3172.2
$ws.Range("H29").Value = 35000
$ws.Range("J29").Value = 35000
$ws.Range("L29").Value = 35000
$ws.Range("N29").Value = -35586
$ws.Range("H31").Value = 18409.936
$ws.Range("J31").Value = 26736.098
$ws.Range("L31").Value = 26736.098
$ws.Range("N31").Value = -27326.098
$ws.Range("H34").Value = 18409.936
$ws.Range("J34").Value = 26736.098
$ws.Range("L34").Value = 26736.098
$ws.Range("N34").Value = -27140.098
$ws.Range("H58").Value = 8720.421
$ws.Range("I58").Value = 12223.2
$ws.Range("J58").Value = 4828.4443
$ws.Range("K58").Value = 12223.2
$ws.Range("L58").Value = 4828.4443
$ws.Range("M58").Value = -12020.2
$ws.Range("N58").Value = -5234.4443
$ws.Range("H113").Value = 2190.8572
$ws.Range("I113").Value = 1172.5
$ws.Range("J113").Value = 2598.2
$ws.Range("K113").Value = 1172.5
$ws.Range("L113").Value = 2598.2
$ws.Range("M113").Value = 997.5
$ws.Range("N113").Value = -6938.2
$ws.Range("H136").Value = 8720.421
$ws.Range("I136").Value = 12223.2
$ws.Range("J136").Value = 4828.4443
$ws.Range("K136").Value = 36669.60000000001
$ws.Range("L136").Value = 14485.3329
$ws.Range("M136").Value = -34119.60000000001
$ws.Range("N136").Value = -19585.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 46362.5
$ws.Range("J37").Value = 46362.5
$ws.Range("L37").Value = 139087.5
$ws.Range("N37").Value = -139311.5
$ws.Range("H137").Value = 2015.1666
$ws.Range("J137").Value = 3516
$ws.Range("L137").Value = 10548
$ws.Range("N137").Value = -20748

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 11247218
$ws.Range("J102").Value = 2168681.8
$ws.Range("L102").Value = 2168681.8
$ws.Range("N102").Value = -2171925.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 50834.777
$ws.Range("I22").Value = 60135.133
$ws.Range("J22").Value = 4333
$ws.Range("K22").Value = 60135.133
$ws.Range("L22").Value = 4333
$ws.Range("M22").Value = -59840.133
$ws.Range("N22").Value = -4923
$ws.Range("H27").Value = 50834.777
$ws.Range("I27").Value = 60135.133
$ws.Range("J27").Value = 4333
$ws.Range("K27").Value = 60135.133
$ws.Range("L27").Value = 4333
$ws.Range("M27").Value = -60028.133
$ws.Range("N27").Value = -4547
$ws.Range("H68").Value = 3245.182
$ws.Range("I68").Value = 3069.7
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 3069.7
$ws.Range("L68").Value = 5000
$ws.Range("M68").Value = -2320.7
$ws.Range("N68").Value = -6498
$ws.Range("H71").Value = 3245.182
$ws.Range("I71").Value = 3069.7
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 15348.5
$ws.Range("L71").Value = 25000
$ws.Range("M71").Value = -11604.5
$ws.Range("N71").Value = -32488
$ws.Range("H100").Value = 147057.42
$ws.Range("I100").Value = 4900.6665
$ws.Range("K100").Value = 4900.6665
$ws.Range("M100").Value = -4359.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 28524.5
$ws.Range("I39").Value = 25000
$ws.Range("J39").Value = 25000
$ws.Range("K39").Value = 25000
$ws.Range("M39").Value = -24587
